# Daily refresh of the cryptos price/volume table (cols D & E), plus a
# rank swap between Bittensor and Fetch.AI (rows 28-29, cols B-E).
# Numeric-looking "Price" strings are written with a leading apostrophe
# so Excel keeps them as literal text instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.215.84"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "2.509.94"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'552.94"
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").Value = "'148.24"
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "2.510.06"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("E10").Value = "  -7.86%  "
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  -7.14%  "
$ws.Range("D13").Value = "'0.358"
$ws.Range("E13").Value = "  -5.66%  "
$ws.Range("D14").Value = "'26.39"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("D15").Value = "2.968.70"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "62.141.62"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'0.0000165"
$ws.Range("E17").Value = "  -6.83%  "
$ws.Range("D18").Value = "2.515.85"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("E19").Value = "  -5.96%  "
$ws.Range("D20").Value = "'7.07"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").Value = "'4.22"
$ws.Range("E21").Value = "  -6.79%  "
$ws.Range("D22").Value = "'324.26"
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'64.17"
$ws.Range("E24").Value = "  -4.31%  "
$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'0.0000105"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("D27").Value = "2.619.52"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "'552.13"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "'1.52"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'8.43"
$ws.Range("E30").Value = "  -7.18%  "
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'7.79"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'0.151"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  -7.03%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  -7.83%  "
$ws.Range("D37").Value = "'4.94"
$ws.Range("E37").Value = "  -8.69%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'0.382"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").Value = "'18.65"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").Value = "'144.45"
$ws.Range("E41").Value = "  -6.76%  "
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -6.59%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'40.69"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "'2.38"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").Value = "'150.43"
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").Value = "'3.60"
$ws.Range("E47").Value = "  -7.65%  "
$ws.Range("D48").Value = "'21.22"
$ws.Range("E48").Value = "  -7.81%  "
$ws.Range("E49").Value = "  -7.81%  "
$ws.Range("D50").Value = "'0.593"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("D51").Value = "'0.0948"
$ws.Range("E51").Value = "  -5.15%  "
